$d = $word.ActiveDocument

# 1. Update the date
$d.Content.Find.Execute("2016.8.16", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2016.8.17", 2)

# 2. Update "today's work result"
$d.Content.Find.Execute("看完了27，28章的内容，敲了书上的代码。", $true, $false, $false, $false, $false,
                         $true, 1, $false, "反复看了哈夫曼编码和树的可视化。做课程设计。", 2)

# 3. Update "issues and improvement"
$d.Content.Find.Execute("感觉有点困难，内容大多都没看懂。只有明天继续看一遍。", $true, $false, $false, $false, $false,
                         $true, 1, $false, "课程设计中如何输入字符集得到哈夫曼树以及如何将树存入文件不会。", 2)

# 4. Update "tomorrow's plan"
$d.Content.Find.Execute("计划将书上的内容在过一遍，准备做课程设计", $true, $false, $false, $false, $false,
                         $true, 1, $false, "继续做课程设计。", 2)

# 5. Remove the now-superfluous trailing empty paragraph that followed the
#    "tomorrow's plan" text in that same table cell. Locate it via the
#    document-level Paragraphs collection (accessing $d.Tables first is
#    known to disturb Paragraphs indexing in this runtime, so it is
#    deliberately avoided here).
$planIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "继续做课程设计。") {
        $planIndex = $i
        break
    }
}

if ($planIndex -ge 1) {
    $nextPara = $d.Paragraphs.Item($planIndex + 1)
    $nextText = $nextPara.Range.Text
    if ($nextText -eq ([string][char]13 + [string][char]7) -or $nextText -eq [string][char]13) {
        $nextPara.Range.Delete()
    }
}

